$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Apply the highlight style (same fill as A1/A2/...) to a few existing cells ---
$ws.Range("A1").Copy()
$ws.Range("A14").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A33").PasteSpecial(-4122)
$ws.Range("A40").PasteSpecial(-4122)
$ws.Range("A41").PasteSpecial(-4122)
$ws.Range("A42").PasteSpecial(-4122)
$ws.Range("A43").PasteSpecial(-4122)
$ws.Range("A44").PasteSpecial(-4122)
$ws.Range("A48").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Add the new "2022" entries, rows 41-63 ---
for ($i = 41; $i -le 63; $i++) {
    $ws.Range("A$i").Value = $i
}

# --- Restore the selection to the new active cell ---
$ws.Range("A42").Select()
